$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with default (unstyled) formatting, used to restore format
# after forcing numeric-looking strings in column D to remain text.
$donor = $ws.Range("C2")

$ws.Range("D2").Value = "30.348.73"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "1.934.50"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("D4").Value = "'1.002"
$donor.Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'251.81"
$donor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  +2.79%  "

$ws.Range("D6").Value = "'0.7229"
$donor.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  +3.12%  "

$ws.Range("D7").Value = "'1.001"
$donor.Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.3291"
$donor.Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("E8").Value = "  +2.19%  "

$ws.Range("D9").Value = "'27.81"
$donor.Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  +6.89%  "

$ws.Range("D10").Value = "'0.07230"
$donor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  +5.91%  "

$ws.Range("D11").Value = "'0.8065"
$donor.Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +2.65%  "

$ws.Range("D12").Value = "'0.08103"
$donor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "1.936.23"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").Value = "'5.460"
$donor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  +2.30%  "

$ws.Range("D15").Value = "'94.69"
$donor.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").Value = "'15.05"
$donor.Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "  +5.02%  "

$ws.Range("D17").Value = "30.349.80"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "'0.000008209"
$donor.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  +5.13%  "

$ws.Range("D19").Value = "'252.88"
$donor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("D20").Value = "'5.820"
$donor.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").Value = "2.189.78"
$ws.Range("E21").Value = "  +0.88%  "

$ws.Range("D22").Value = "'1.001"
$donor.Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'6.942"
$donor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +2.07%  "

$ws.Range("D25").Value = "'9.723"
$donor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  +2.15%  "

$ws.Range("D26").Value = "'166.08"
$donor.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = "  +3.75%  "

$ws.Range("D27").Value = "'2.345"
$donor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  +5.80%  "

$ws.Range("D28").Value = "'19.27"
$donor.Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  +3.28%  "

$ws.Range("D29").Value = "'0.1297"
$donor.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").Value = "'1.355"
$donor.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").Value = "'1.549"
$donor.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("D32").Value = "'4.436"
$donor.Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").Value = "'4.196"
$donor.Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("D34").Value = "'0.05244"
$donor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = "  +4.51%  "

$ws.Range("D35").Value = "'1.262"
$donor.Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("E35").Value = "  +6.75%  "

$ws.Range("D36").Value = "'0.7503"
$donor.Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("D37").Value = "'2.771"
$donor.Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("D38").Value = "'0.01966"
$donor.Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").Value = "  +2.57%  "

$ws.Range("D39").Value = "'2.804"
$donor.Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").Value = "'79.28"
$donor.Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").Value = "'6.442"
$donor.Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("D42").Value = "'0.4534"
$donor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  +3.02%  "

$ws.Range("D43").Value = "'2.033"
$donor.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("D44").Value = "'0.8460"
$donor.Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  +1.61%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").Value = "'101.95"
$donor.Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").Value = "'9.808"
$donor.Copy()
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("D48").Value = "'7.439"
$donor.Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  +3.74%  "

$ws.Range("D49").Value = "'36.74"
$donor.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("D50").Value = "'0.4193"
$donor.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  +3.59%  "

$ws.Range("D51").Value = "'0.06050"
$donor.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  +2.29%  "
